# "correction 4 avec cours de M"
# Updates Sheet1 content: new/extended explanations on Burndown/KPI, the
# "Developers" definition-of-done paragraph, the user-story/debt paragraph,
# the scrum-team composition paragraph, and the scrum-history date cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# ---------------------------------------------------------------------
# A4 : "Burndown chart" -> "Burndup/Burndown chart, KPI"
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Burndup/Burndown chart, KPI"

# ---------------------------------------------------------------------
# B4 : burndown/KPI explanation, extended with new KPI bullet points
# ---------------------------------------------------------------------
$ws.Range("B4").Value = (
    "How much work remains till the end of the sprint, shows the evolution of remaining effort against time. " + $nl +
    "NEW : KPI > cumulative flow diagram : reduire la colonne ""en cours""." + $nl +
    "NEW : KPI > vélocité :  des story points (mesure sur 5 sprints), pour savoir nbre pts/sprint." + $nl +
    "UP : Outil non précis et ne prévoit pas." + $nl +
    "DOWN : analyse de tendance et prévision." + $nl +
    "Responsable : developers."
)

# Row 4 grows taller to fit the extra lines.
$ws.Rows.Item(4).RowHeight = 107.4

# ---------------------------------------------------------------------
# B26 : "Developers (...)" rich-text paragraph, extra clause inserted in
# the last (plain) run.
# ---------------------------------------------------------------------
$b26Parts = @(
    @{ Text = 'Developers (responsable du suivi du travail restant), la définition de "done" '; Bold = $false; Underline = $false },
    @{ Text = 'créée par les developers'; Bold = $false; Underline = $true },
    @{ Text = ' peut possiblement '; Bold = $false; Underline = $false },
    @{ Text = "s'adapter"; Bold = $false; Underline = $true },
    @{ Text = ' à chaque RETRO '; Bold = $false; Underline = $false },
    @{ Text = 'par la scrum team'; Bold = $false; Underline = $true },
    @{ Text = ", n'aide pas à calculer la vélocité." + $nl + "DONE est l'objectif, "; Bold = $false; Underline = $false },
    @{ Text = 'pas la réduction de la dette technique par un sprint spécial (hardening)'; Bold = $true; Underline = $false },
    @{ Text = '.' + $nl + 'Doit respecter : conventions & standards, same definition for other teams working on the same product, pas de travail supplémentaire à faire.' + $nl + 'ensures artifact transparency, is used to acces, guides the developers.'; Bold = $false; Underline = $false }
)

$b26Cell = $ws.Range("B26")
$full = ""
foreach ($p in $b26Parts) { $full += $p.Text }
$b26Cell.Value = $full

$pos = 1
foreach ($p in $b26Parts) {
    $len = $p.Text.Length
    if ($p.Bold -or $p.Underline) {
        $chars = $b26Cell.Characters($pos, $len)
        if ($p.Bold) { $chars.Font.Bold = $true }
        if ($p.Underline) { $chars.Font.Underline = $true }
    }
    $pos += $len
}

# ---------------------------------------------------------------------
# B44 : user-story / velocity rich-text paragraph, new "Dette technique"
# runs appended at the end.
# ---------------------------------------------------------------------
$b44Parts = @(
    @{ Text = "user-story, velocity, definition of ready >> tous ces termes n'existent pas dans le guide scrum." + $nl + "La vélocité est une mesure de la "; Bold = $false; Underline = $false },
    @{ Text = 'maturité'; Bold = $true; Underline = $false },
    @{ Text = " de l'équipe, pas un indice de réussite." + $nl + "Composition des équipes de dev en fonction de la "; Bold = $false; Underline = $false },
    @{ Text = '"self-organization"'; Bold = $true; Underline = $false },
    @{ Text = '.' + $nl + 'Developers : '; Bold = $false; Underline = $false },
    @{ Text = 'Adjustments'; Bold = $true; Underline = $false },
    @{ Text = ' to its engineering practices whenever needed.' + $nl; Bold = $false; Underline = $false },
    @{ Text = 'User story'; Bold = $true; Underline = $false },
    @{ Text = ' : se rédige comme suit : en tant que, je veux, afin de. ' + $nl + 'Méthode '; Bold = $false; Underline = $false },
    @{ Text = 'Moscow'; Bold = $true; Underline = $false },
    @{ Text = ' pour '; Bold = $false; Underline = $false },
    @{ Text = 'prioriser'; Bold = $false; Underline = $true },
    @{ Text = ' selon les points accordés aux Should et Could.' + $nl + 'Planning Poker : story points estimés avec la business value et la complexité de dév.' + $nl + '(points des dev : prendre le max et le min, si trop grand écart, réévaluer sinon la moyenne).' + $nl; Bold = $false; Underline = $false },
    @{ Text = 'Dette technique : '; Bold = $true; Underline = $false },
    @{ Text = 'non respect de la conception voulue ou non induisant des coûts supplémentaires vus comme des interêts.' + $nl; Bold = $false; Underline = $false }
)

$b44Cell = $ws.Range("B44")
$full = ""
foreach ($p in $b44Parts) { $full += $p.Text }
$b44Cell.Value = $full

$pos = 1
foreach ($p in $b44Parts) {
    $len = $p.Text.Length
    if ($p.Bold -or $p.Underline) {
        $chars = $b44Cell.Characters($pos, $len)
        if ($p.Bold) { $chars.Font.Bold = $true }
        if ($p.Underline) { $chars.Font.Underline = $true }
    }
    $pos += $len
}

# Row 44 grows taller to fit the extra "Dette technique" paragraph.
$ws.Rows.Item(44).RowHeight = 156.6

# ---------------------------------------------------------------------
# B30 : scrum-team composition paragraph, extra closing sentence added
# ---------------------------------------------------------------------
$ws.Range("B30").Value = (
    "Scrum master, product owner, developers.  " + $nl +
    "Différentes équipes scrum d'un même produit peuvent NE PAS avoir des sprints de longueur différente (does not require) et ne pas avoir de sprint ""alignés""." + $nl +
    "Elle fait le SPRINT PLANNING, responsible for crafting the sprint goal during it. Elle doit être autosuffisante." + $nl +
    "Qualités importantes : flexibility, creativity, productivity." + $nl +
    "Should have all competencies, should choose how best to accomplish their work." + $nl +
    "If a 2nd scrum team is added, productivity is likely to decrease." + $nl +
    "BEFORE : self-organizing, NOW : self-managing."
)

# ---------------------------------------------------------------------
# C14 : scrum history sentence, updated dates/guide count
# ---------------------------------------------------------------------
$ws.Range("C14").Value = "Le mot scrum est apparu en 1995, l'agilité en 2001, le 1er scrum guide en 2011, le 2ème en 2020."

# ---------------------------------------------------------------------
# Move the saved cursor/selection from C34 to C14 (matches the author's
# last edited cell).
# ---------------------------------------------------------------------
$ws.Range("C14").Select()

Write-Host "Edit complete"
